# Fruta / hortaliza, semanal
# Insert a new weekly observation row for Limón (Agrícola del Norte S.A. de
# Arica) at row 131, pushing the existing rows 131-221 down to 132-222.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 131:221 down by one to make room for the new record.
$ws.Rows("131:131").Insert()

# Populate the newly inserted row 131 with the new weekly record.
$ws.Range("A131").Value = 1
$ws.Range("B131").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C131").Value = "Arica y Parinacota"
$ws.Range("D131").Value = 44651
$ws.Range("E131").Value = 15
$ws.Range("F131").Value = "Fruta"
$ws.Range("G131").Value = 100102
$ws.Range("H131").Value = "Cítricos"
$ws.Range("I131").Value = 100102003
$ws.Range("J131").Value = "Limón"
$ws.Range("K131").Value = "Sin especificar"
$ws.Range("L131").Value = "2a amarillo"
$ws.Range("M131").Value = 270
$ws.Range("N131").Value = 30000
$ws.Range("O131").Value = 31000
$ws.Range("P131").Value = 30500
$ws.Range("Q131").Value = "$/caja 20 kilos"
$ws.Range("R131").Value = "Región de Coquimbo"
$ws.Range("S131").Value = 1525
$ws.Range("T131").Value = 20
